$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'51.872.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
$ws.Range("D3").Value = "'2.782.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'358.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.43%  "

# Row 6
$ws.Range("D6").Value = "'109.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.85%  "

# Row 7
$ws.Range("D7").Value = "'0.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.40%  "

# Row 8
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.64%  "

# Row 10
$ws.Range("D10").Value = "'40.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.85%  "

# Row 11
$ws.Range("D11").Value = "'0.0853"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "

# Row 12
$ws.Range("E12").Value = "  +0.72%  "

# Row 13
$ws.Range("D13").Value = "'19.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "

# Row 14
$ws.Range("E14").Value = "  -1.64%  "

# Row 15
$ws.Range("D15").Value = "'3.217.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.03%  "

# Row 16
$ws.Range("D16").Value = "'2.800.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.40%  "

# Row 17
$ws.Range("E17").Value = "  +3.75%  "

# Row 18
$ws.Range("D18").Value = "'51.804.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "

# Row 19
$ws.Range("E19").Value = "  +0.85%  "

# Row 20
$ws.Range("E20").Value = "  -1.20%  "

# Row 21
$ws.Range("D21").Value = "'13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "

# Row 22
$ws.Range("E22").Value = "  -1.75%  "

# Row 23
$ws.Range("D23").Value = "'274.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "

# Row 24
$ws.Range("D24").Value = "'70.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "

# Row 25
$ws.Range("D25").Value = "'2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "

# Row 26
$ws.Range("E26").Value = "  -0.32%  "

# Row 27
$ws.Range("E27").Value = "  -0.08%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.82%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'10.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("E30").Value = "  +3.83%  "

# Row 31
$ws.Range("E31").Value = "  +5.15%  "

# Row 32
$ws.Range("D32").Value = "'51.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "

# Row 33
$ws.Range("D33").Value = "'33.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.69%  "

# Row 34
$ws.Range("E34").Value = "  -2.40%  "

# Row 35
$ws.Range("E35").Value = "  +2.29%  "

# Row 36
$ws.Range("E36").Value = "  +6.86%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("E38").Value = "  +1.23%  "

# Row 39
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("E40").Value = "  -3.71%  "

# Row 41
$ws.Range("E41").Value = "  +0.92%  "

# Row 42
$ws.Range("E42").Value = "  -1.19%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'121.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.42%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.42%  "

# Row 45
$ws.Range("D45").Value = "'22.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.19%  "

# Row 46
$ws.Range("D46").Value = "'2.069.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "

# Row 47
$ws.Range("E47").Value = "  -2.60%  "

# Row 48
$ws.Range("D48").Value = "'2.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.11%  "

# Row 49
$ws.Range("D49").Value = "'5.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.08%  "

# Row 50
$ws.Range("E50").Value = "  -0.68%  "

# Row 51
$ws.Range("D51").Value = "'8.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
